$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.360.27"
$ws.Range('E2').Value = '  -0.34%  '

$ws.Range('D3').Value = "'1.844.37"
$ws.Range('E3').Value = '  -0.43%  '

$ws.Range('D4').Value = "'0.9994"
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "'238.64"
$ws.Range('E5').Value = '  -1.35%  '

$ws.Range('D6').Value = "'0.6280"
$ws.Range('E6').Value = '  -0.48%  '

$ws.Range('D7').Value = "'1.000"
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = "'0.07547"
$ws.Range('E8').Value = '  -0.80%  '

$ws.Range('D9').Value = "'0.2942"
$ws.Range('E9').Value = '  -1.16%  '

$ws.Range('D10').Value = "'24.59"
$ws.Range('E10').Value = '  +0.31%  '

$ws.Range('D11').Value = "'0.07696"
$ws.Range('E11').Value = '  -0.24%  '

$ws.Range('D12').Value = "'1.845.60"
$ws.Range('E12').Value = '  -3.98%  '

$ws.Range('D13').Value = "'4.978"
$ws.Range('E13').Value = '  -0.61%  '

$ws.Range('D14').Value = "'0.6787"
$ws.Range('E14').Value = '  -1.29%  '

$ws.Range('D15').Value = "'0.00001023"
$ws.Range('E15').Value = '  +3.21%  '

$ws.Range('D16').Value = "'83.09"

$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = "'2.106.50"
$ws.Range('E17').Value = '  -4.76%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = "'6.124"
$ws.Range('E18').Value = '  -0.96%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = "'29.406.87"
$ws.Range('E19').Value = '  -0.66%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'228.13"
$ws.Range('E20').Value = '  -2.17%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'12.43"
$ws.Range('E21').Value = '  -0.80%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'1.000"
$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = "'7.494"
$ws.Range('E23').Value = '  -2.88%  '

$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = "'1.001"
$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'156.95"
$ws.Range('E25').Value = '  +1.39%  '

$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = "'0.1390"
$ws.Range('E26').Value = '  -0.69%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = "'8.357"
$ws.Range('E27').Value = '  -1.37%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'17.63"
$ws.Range('E28').Value = '  -0.19%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = "'1.461"
$ws.Range('E29').Value = '  -0.88%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = "'1.263"
$ws.Range('E30').Value = '  +0.38%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.05627"
$ws.Range('E31').Value = '  -2.81%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'4.117"
$ws.Range('E32').Value = '  -0.39%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = "'4.031"
$ws.Range('E33').Value = '  +0.30%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = "'1.837"
$ws.Range('E34').Value = '  -2.95%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'1.154"
$ws.Range('E35').Value = '  -1.12%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'0.7162"
$ws.Range('E36').Value = '  -0.90%  '

$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = "'2.594"
$ws.Range('E37').Value = '  +0.30%  '

$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = "'1.245.73"
$ws.Range('E38').Value = '  -0.41%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.01806"
$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = "'2.773"
$ws.Range('E40').Value = '  -0.72%  '

$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'6.203"
$ws.Range('E41').Value = '  +1.75%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'0.9031"
$ws.Range('E42').Value = '  -0.59%  '

$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = "'1.000"
$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'101.65"
$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'66.24"
$ws.Range('E45').Value = '  -2.26%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = "'0.00000000122"
$ws.Range('E46').Value = '  +3.89%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = "'7.093"
$ws.Range('E47').Value = '  -3.04%  '

$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').Value = "'0.3997"
$ws.Range('E48').Value = '  -1.03%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'9.008"
$ws.Range('E49').Value = '  -1.56%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'1.684"
$ws.Range('E50').Value = '  -1.65%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.1117"
$ws.Range('E51').Value = '  -0.48%  '

